$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Materialize explicit (but empty) cells for the new D:K columns on rows 1-14,
# matching the OOXML produced by Excel when a whole block is touched.
$ws.Range("D1:K14").Font.Bold = $false

# Populate the new "remark" column (E) and the scattered extra columns
# (H / I / K) with their values.
$ws.Range("E1").Value  = "备注"
$ws.Range("E2").Value  = "空字符串"
$ws.Range("E3").Value  = "备注1"
$ws.Range("E4").Value  = "备注2"
$ws.Range("H4").Value  = "是非得失"
$ws.Range("E5").Value  = "备注3"
$ws.Range("E6").Value  = "备注4"
$ws.Range("E7").Value  = "备注5"
$ws.Range("E8").Value  = "备注6"
$ws.Range("H8").Value  = "是非得失"
$ws.Range("K8").Value  = "水电费"
$ws.Range("I9").Value  = "是否"
$ws.Range("H14").Value = "是非得失"
